# [PBL] Neue Zusätzliches PBI für Sprint 4.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "PBL2" to "PBL"
$ws.Name = "PBL"

# PBI #12 ("Als PO möchte ich für jedes PBI einen Forecast ...") is now Done
$ws.Cells.Item(8, 5).Value = "Done"

# Insert a new row for the new PBI #13 right above the old row 9
$ws.Rows.Item(9).Insert()

$ws.Cells.Item(9, 1).Value = 13
$ws.Cells.Item(9, 2).Value = "Als PO möchte ich Releases durch die Angabe der ID eines PBIs definieren können. Für diese Releases werden die Forecasts berechnet, bis zu welchen Sprint diese fertiggestellt sind."
$ws.Cells.Item(9, 3).Value = "Akzeptanzkriteren:`n- Die Releases werden in einer Taballe auf einem neuen Tab dargestellt.`n- Die Releases können editiert und abgespeichert werden.`n- Es wird angezeigt, wenn es kein PBI mit dieser ID gibt."
$ws.Cells.Item(9, 4).Value = 3
$ws.Cells.Item(9, 5).Value = "Todo"
$ws.Cells.Item(9, 6).Value = "Sprint 4"

# Match the wrap-text style used by the other Title/Summary cells
$ws.Cells.Item(9, 2).WrapText = $true
$ws.Cells.Item(9, 3).WrapText = $true

# Row height for the new row
$ws.Rows.Item(9).RowHeight = 120

# Update selection to match the new edit location
$ws.Range("B9").Select()

# Page setup as captured in the saved workbook
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
